$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UiPathForm")

# Add a new row of shopping-list data (row 11). Copy the previous data
# row (row 10) down to row 11 first so the new row inherits the exact
# same cell formatting/styles as the rest of the table.
$ws.Range("A10:M10").Copy()
$ws.Range("A11:M11").PasteSpecial(-4104)   # xlPasteAll
$excel.CutCopyMode = 0

# Plain text values can simply be assigned - Excel keeps them as text.
$ws.Range("F11").Value = "Green Tea"
$ws.Range("G11").Value = "1L"
$ws.Range("I11").Value = "Bread"
$ws.Range("J11").Value = "1Pkt"

# H11 and K11 need the digit-only text values "2" and "3" respectively.
# Assigning those directly would make Excel re-interpret them as numbers,
# so instead copy the already-existing text cells that hold those exact
# values and paste just the values in, which preserves their text type.
$ws.Range("H8").Copy()
$ws.Range("H11").PasteSpecial(-4163)       # xlPasteValues ("2")
$excel.CutCopyMode = 0

$ws.Range("E7").Copy()
$ws.Range("K11").PasteSpecial(-4163)       # xlPasteValues ("3")
$excel.CutCopyMode = 0

# A11, B11, C11, D11, E11, L11, M11 already match row 10's values, and
# were brought over correctly by the row copy above, so no further
# changes are required for them.
